$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "pythonCode": turn on wrap-text for the search-code rows (A4:A8),
# which introduces a new cell style (font2/fill2/left-align + wrapText).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pythonCode")
$ws2.Range("A4:A8").WrapText = $true

# ---------------------------------------------------------------------------
# Sheet "Sheet1": insert a brand-new row 7 (copies formatting from row 6,
# so columns A-D inherit style index 1 automatically), then populate it.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(7).Insert()

$ws1.Range("A7").Value = "Numpysdet84"
$ws1.Range("B7").Value = "sdet84batch"
$ws1.Range("C7").Value = "You are logged in"

$searchCode = "def search(input_list, num):`nif(num in input_list):`nprint(`"Element Found`")`n\b`n\b`nelse:`nprint(`"Not Found`")`n\b`n\b`n\b`n\b`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$ws1.Range("D7").Value = $searchCode
$ws1.Range("D7").WrapText = $true

$ws1.Range("E7").Value = "Element Found"

# Inserting/auto-wrapping can push the row height up to fit the long text;
# restore the same fixed height used by every other data row.
$ws1.Rows.Item(7).RowHeight = 15

# Widen column D so the new multi-line code sample is readable.
$ws1.Columns.Item(4).ColumnWidth = 128.94401041666666

# ---------------------------------------------------------------------------
# Update the remembered selections on both sheets.
# ---------------------------------------------------------------------------
$ws2.Range("A4").Select() | Out-Null
$ws1.Range("D7").Select() | Out-Null

# Sheet1 was (and stays) the active/visible tab.
$ws1.Activate() | Out-Null
